# Updates the cryptos worksheet: Price (column D) and Volume(1h) (column E)
# values for rows 2-51, matching the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to store the value as text (matching the original
    # inlineStr/shared-string representation) instead of letting Excel
    # auto-convert numeric-looking strings (e.g. "0.635") into numbers.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "37.761.67"
Set-TextValue $ws.Range("E2") "  +0.12%  "
Set-TextValue $ws.Range("D3") "2.087.25"
Set-TextValue $ws.Range("E3") "  +0.58%  "
Set-TextValue $ws.Range("E4") "  -0.09%  "
Set-TextValue $ws.Range("D5") "233.69"
Set-TextValue $ws.Range("E5") "  +0.07%  "
Set-TextValue $ws.Range("D6") "0.635"
Set-TextValue $ws.Range("E6") "  +1.94%  "
Set-TextValue $ws.Range("E7") "  -0.01%  "
Set-TextValue $ws.Range("D8") "58.26"
Set-TextValue $ws.Range("E8") "  +0.38%  "
Set-TextValue $ws.Range("D9") "0.393"
Set-TextValue $ws.Range("E9") "  +0.65%  "
Set-TextValue $ws.Range("D10") "0.0782"
Set-TextValue $ws.Range("E10") "  +0.11%  "
Set-TextValue $ws.Range("E11") "  +3.00%  "
Set-TextValue $ws.Range("D12") "15.22"
Set-TextValue $ws.Range("E12") "  +3.03%  "
Set-TextValue $ws.Range("D13") "2.392.77"
Set-TextValue $ws.Range("E13") "  +0.37%  "
Set-TextValue $ws.Range("D14") "21.14"
Set-TextValue $ws.Range("E14") "  +1.43%  "
Set-TextValue $ws.Range("D15") "0.780"
Set-TextValue $ws.Range("E15") "  +1.22%  "
Set-TextValue $ws.Range("D16") "5.35"
Set-TextValue $ws.Range("E16") "  +1.18%  "
Set-TextValue $ws.Range("D17") "2.082.38"
Set-TextValue $ws.Range("E17") "  +1.23%  "
Set-TextValue $ws.Range("D18") "37.754.89"
Set-TextValue $ws.Range("E18") "  +0.16%  "
Set-TextValue $ws.Range("D19") "6.10"
Set-TextValue $ws.Range("E19") "  -1.05%  "
Set-TextValue $ws.Range("D20") "71.01"
Set-TextValue $ws.Range("E20") "  -0.15%  "
Set-TextValue $ws.Range("D21") "0.0₃0837"
Set-TextValue $ws.Range("E21") "  +0.69%  "
Set-TextValue $ws.Range("D22") "229.47"
Set-TextValue $ws.Range("E22") "  +0.82%  "
Set-TextValue $ws.Range("E23") "  -0.04%  "
Set-TextValue $ws.Range("E24") "  -0.77%  "
Set-TextValue $ws.Range("E25") "  +0.03%  "
Set-TextValue $ws.Range("D26") "9.74"
Set-TextValue $ws.Range("E26") "  +8.59%  "
Set-TextValue $ws.Range("D27") "170.88"
Set-TextValue $ws.Range("E27") "  +0.85%  "
Set-TextValue $ws.Range("D28") "0.133"
Set-TextValue $ws.Range("E28") "  -3.90%  "
Set-TextValue $ws.Range("D29") "19.54"
Set-TextValue $ws.Range("E29") "  +0.59%  "
Set-TextValue $ws.Range("E30") "  -0.06%  "
Set-TextValue $ws.Range("D31") "0.123"
Set-TextValue $ws.Range("E31") "  +1.38%  "
Set-TextValue $ws.Range("D32") "4.69"
Set-TextValue $ws.Range("E32") "  +0.46%  "
Set-TextValue $ws.Range("E33") "  +1.09%  "
Set-TextValue $ws.Range("D34") "4.60"
Set-TextValue $ws.Range("E34") "  -1.17%  "
Set-TextValue $ws.Range("D35") "2.50"
Set-TextValue $ws.Range("E35") "  +1.82%  "
Set-TextValue $ws.Range("E36") "  -0.24%  "
Set-TextValue $ws.Range("D37") "3.33"
Set-TextValue $ws.Range("E37") "  -1.61%  "
Set-TextValue $ws.Range("E38") "  -0.25%  "
Set-TextValue $ws.Range("D39") "5.38"
Set-TextValue $ws.Range("E39") "  +0.66%  "
Set-TextValue $ws.Range("D40") "0.0236"
Set-TextValue $ws.Range("E40") "  +9.75%  "
Set-TextValue $ws.Range("D41") "101.35"
Set-TextValue $ws.Range("E41") "  +3.52%  "
Set-TextValue $ws.Range("D42") "0.0964"
Set-TextValue $ws.Range("E42") "  -1.16%  "
Set-TextValue $ws.Range("E43") "  +1.24%  "
Set-TextValue $ws.Range("D44") "1.20"
Set-TextValue $ws.Range("E44") "  +3.41%  "
Set-TextValue $ws.Range("D45") "16.70"
Set-TextValue $ws.Range("E45") "  +1.43%  "
Set-TextValue $ws.Range("D46") "1.451.82"
Set-TextValue $ws.Range("E46") "  -0.05%  "
Set-TextValue $ws.Range("D47") "4.11"
Set-TextValue $ws.Range("E47") "  -3.04%  "
Set-TextValue $ws.Range("D48") "1.06"
Set-TextValue $ws.Range("E48") "  -0.24%  "
Set-TextValue $ws.Range("D49") "7.20"
Set-TextValue $ws.Range("E49") "  -2.28%  "
Set-TextValue $ws.Range("E50") "  -1.96%  "
Set-TextValue $ws.Range("D51") "2.276.75"
Set-TextValue $ws.Range("E51") "  +0.35%  "
